# Applies the edit described by the diff:
#  1) Column F (trial_total) is decremented by 81 for every data row (2..42).
#  2) The "stimulus block" (columns L..V: stimulus, conceptual, perceptual,
#     typicality, n, p_typicality, p_conceptual, p_perceptual, r_typicality,
#     r_conceptual, r_perceptual) is reshuffled among a subset of rows
#     (a row receives the block that used to live on another row).
#  3) Row 12's catch-trial stimulus filename is renamed independently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shift trial_total (column F) down by 81 for rows 2..42 ---
for ($r = 2; $r -le 42; $r++) {
    $cur = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 6).Value = $cur - 81
}

# --- 2) Permute the stimulus block (columns L..V = 12..22) ---
# Maps destination row -> source row (the row whose L..V values move there).
$rowMap = @{
    4  = 21
    7  = 4
    9  = 7
    10 = 24
    11 = 32
    14 = 9
    16 = 35
    20 = 10
    21 = 11
    23 = 26
    24 = 14
    26 = 16
    29 = 39
    30 = 20
    32 = 41
    35 = 23
    39 = 29
    41 = 30
}

$cols = 12..22

# Snapshot the "before" values for every row involved, so overwrites don't
# clobber data that another destination still needs to read (the mapping
# contains cycles).
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# Apply the snapshot values to their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c]
    }
}

# --- 3) Rename the catch-trial stimulus on row 12 ---
$ws.Cells.Item(12, 12).Value = "stimuli/catch_25.jpg"
